# From v1.2 to v1.2.1
# Reorders the "TC2" / "TC3" / "TC4" test-case bodies in the UC011 sheet:
#   - TC2's second step becomes the old TC4 content (ordenar pelo nome do servidor)
#   - TC3's second step becomes the old TC2 content (busca/filtro)
#   - TC4's second step becomes the old TC3 content (cancelamento de diária)
#   - TC5 is untouched
# The "TCx" labels themselves (B15/B23/B31/B39) stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ordenarAcao = "Chefe Clica para ordenar pelo nome do servidor."
$ordenarResultado = "SYSTEM Visualiza os registros de solicitações de diária ordenado pelo nome do servidor."

$buscaAcao = "Chefe Indica alguns parâmetros específicos para a busca; Informa o nome do beneficiário; Filtra a listagem de solicitações."
$buscaResultado = "SYSTEM Exibe uma nova listagem de solicitações, de acordo com os filtros informados pelo usuário."

$cancelamentoAcao = "Chefe Clica para realizar o cancelamento de uma diária."
$cancelamentoResultado = "SYSTEM Verifica que a solicitação está em situação SOLICITADA; Exibe mensagem de confirmação (MSG987 - Cancelar solicitação de diária) para o usuário (que deve confirmar); Cancela a diária, mudando sua situação para CANCELADA (ver diagrama de estados da diária)."

# TC2 block (rows 15-20): second step goes from "busca" to "ordenar"
$ws.Range("B20").Value = $ordenarAcao
$ws.Range("D20").Value = $ordenarResultado

# TC3 block (rows 23-28): second step goes from "cancelamento" to "busca"
$ws.Range("B28").Value = $buscaAcao
$ws.Range("D28").Value = $buscaResultado

# TC4 block (rows 31-36): second step goes from "ordenar" to "cancelamento"
$ws.Range("B36").Value = $cancelamentoAcao
$ws.Range("D36").Value = $cancelamentoResultado
